$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.877.08"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.279.06"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.30"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "78.91"
$ws.Range("E7").Value = "  +8.43%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.644"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.45"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0968"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "2.619.15"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.21"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.869"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "2.280.11"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "42.783.15"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "0.0₃0997"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.27"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.17"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.19"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.79"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.34"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.35"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.61"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0857"
$ws.Range("E32").Value = "  +5.15%  "
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.30"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.75"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.91"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "114.67"
$ws.Range("E42").Value = "  +17.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.210"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.48"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.61"
$ws.Range("E47").Value = "  -7.80%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.29"
$ws.Range("E51").Value = "  -0.72%  "
